$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (S1_TC_T4): rename test description
$ws.Range("B5").Value = "Create comment unsupported html format"

# Row 11 (S1_TC_T10): append stored comment id to validation
$ws.Range("J11").Value = "status=200||targetType=TRRecord||targetId=456539938WOS1||id=(S1_TC_T1_comments.id)"

# Row 12 (S1_TC_T11): append stored comment id to validation
$ws.Range("J12").Value = "status=200||targetType=TRRecord||targetId=456539938WOS1||id=(S1_TC_T2_comments.id)"

# Row 13 (S1_TC_T12): append stored comment id to validation
$ws.Range("J13").Value = "status=200||targetType=TRRecord||targetId=456539938WOS1||id=(S1_TC_T3_comments.id)"

# Row 8 (S1_TC_T7): append stored comment id to validation
$ws.Range("J8").Value = "status=200||content=hi||targetType=TRRecord||targetId=456539938WOS1||id=(S1_TC_T1_comments.id)"

# Row 9 (S1_TC_T8): append stored comment id to validation
$ws.Range("J9").Value = "status=200||hasAppreciated=UP||targetType=Comment||appreciateCount=1||targetId=(S1_TC_T1_comments.id)"

# Row 10 (S1_TC_T9): append stored comment id to validation
$ws.Range("J10").Value = "status=200||hasAppreciated=DOWN||targetType=Comment||appreciateCount=0||targetId=(S1_TC_T1_comments.id)"

# Row 3 (S1_TC_T2): switch header user and validation user to SYS_USER2
$ws.Range("F3").Value = "X-1P-User=(SYS_USER2)||Content-Type=application/json"
$ws.Range("J3").Value = "status=200||comments.userId=(SYS_USER2)||comments.targetType=TRRecord||comments.targetId=456539938WOS1||comments.content=mohana.yalamarthi@thomsonreuters.com"

# Row 12 (S1_TC_T11) header also switches to SYS_USER2 (reuses shared string)
$ws.Range("F12").Value = "X-1P-User=(SYS_USER2)||Content-Type=application/json"

# Adjust column J width (bestFit width grew because of longer strings)
$ws.Columns.Item(10).ColumnWidth = 108.7109375

# Adjust sheet view: scroll position and active selection
$ws.Range("F12").Select()
$excel.ActiveWindow.ScrollColumn = 5
